$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 147.58
$ws.Range("I15").Value = 147.58
$ws.Range("K15").Value = 442.74
$ws.Range("M15").Value = -273.74
# Row 70
$ws.Range("H70").Value = 1220
$ws.Range("I70").Value = 1050
$ws.Range("J70").Value = 1262.5
$ws.Range("K70").Value = 3150
$ws.Range("L70").Value = 3787.5
$ws.Range("M70").Value = -2880
$ws.Range("N70").Value = -4327.5
# Row 73
$ws.Range("H73").Value = 1220
$ws.Range("I73").Value = 1050
$ws.Range("J73").Value = 1262.5
$ws.Range("K73").Value = 3150
$ws.Range("L73").Value = 3787.5
$ws.Range("M73").Value = -2214
$ws.Range("N73").Value = -5659.5
# Row 137
$ws.Range("H137").Value = 23811738
$ws.Range("I137").Value = 1269.6364
$ws.Range("J137").Value = 50003256
$ws.Range("K137").Value = 3808.9092
$ws.Range("L137").Value = 150009768
$ws.Range("M137").Value = -1258.9092
$ws.Range("N137").Value = -150014868

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14749.857
$ws.Range("I32").Value = 9536.316999999999
$ws.Range("J32").Value = 48838.383
$ws.Range("K32").Value = 9536.316999999999
$ws.Range("L32").Value = 48838.383
$ws.Range("M32").Value = -9249.316999999999
$ws.Range("N32").Value = -49412.383
# Row 45
$ws.Range("H45").Value = 1129.4546
$ws.Range("I45").Value = 969.35
$ws.Range("J45").Value = 1375.7693
$ws.Range("K45").Value = 969.35
$ws.Range("L45").Value = 1375.7693
$ws.Range("M45").Value = -592.35
$ws.Range("N45").Value = -2129.7693
# Row 61
$ws.Range("H61").Value = 3203.3635
$ws.Range("I61").Value = 3164.8235
$ws.Range("J61").Value = 3244.3125
$ws.Range("K61").Value = 3164.8235
$ws.Range("L61").Value = 3244.3125
$ws.Range("M61").Value = -2952.8235
$ws.Range("N61").Value = -3668.3125
# Row 63
$ws.Range("H63").Value = 2304.75
$ws.Range("I63").Value = 1643.7931
$ws.Range("J63").Value = 4047.2727
$ws.Range("K63").Value = 1643.7931
$ws.Range("L63").Value = 4047.2727
$ws.Range("M63").Value = -957.7931000000001
$ws.Range("N63").Value = -5419.2727
# Row 66
$ws.Range("H66").Value = 2304.75
$ws.Range("I66").Value = 1643.7931
$ws.Range("J66").Value = 4047.2727
$ws.Range("K66").Value = 8218.9655
$ws.Range("L66").Value = 20236.3635
$ws.Range("M66").Value = -4786.9655
$ws.Range("N66").Value = -27100.3635
# Row 110
$ws.Range("H110").Value = 1314.7693
$ws.Range("I110").Value = 1210.2
$ws.Range("J110").Value = 1663.3334
$ws.Range("K110").Value = 1210.2
$ws.Range("L110").Value = 1663.3334
$ws.Range("M110").Value = 834.8
$ws.Range("N110").Value = -5753.3334
# Row 122
$ws.Range("H122").Value = 2060.4707
$ws.Range("I122").Value = 1795.3334
$ws.Range("J122").Value = 2358.75
$ws.Range("K122").Value = 5386.0002
$ws.Range("L122").Value = 7076.25
$ws.Range("M122").Value = -2936.0002
$ws.Range("N122").Value = -11976.25
# Row 132
$ws.Range("H132").Value = 32668658
$ws.Range("I132").Value = 71667860
$ws.Range("J132").Value = 169324
$ws.Range("K132").Value = 215003580
$ws.Range("L132").Value = 507972
$ws.Range("M132").Value = -215001050
$ws.Range("N132").Value = -513032
# Row 136
$ws.Range("H136").Value = 3203.3635
$ws.Range("I136").Value = 3164.8235
$ws.Range("J136").Value = 3244.3125
$ws.Range("K136").Value = 9494.470499999999
$ws.Range("L136").Value = 9732.9375
$ws.Range("M136").Value = -6944.470499999999
$ws.Range("N136").Value = -14832.9375

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1841.1154
$ws.Range("I20").Value = 1906.4117
$ws.Range("J20").Value = 1717.7778
$ws.Range("K20").Value = 1906.4117
$ws.Range("L20").Value = 1717.7778
$ws.Range("M20").Value = -1659.4117
$ws.Range("N20").Value = -2211.7778
# Row 134
$ws.Range("H134").Value = 129380.445
$ws.Range("I134").Value = 145178
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 435534
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -432999
$ws.Range("N134").Value = -14070

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 10615.2
$ws.Range("J50").Value = 10615.2
$ws.Range("L50").Value = 10615.2
$ws.Range("N50").Value = -11865.2
# Row 58
$ws.Range("H58").Value = 3198.3936
$ws.Range("I58").Value = 1374.25
$ws.Range("J58").Value = 4381.6216
$ws.Range("K58").Value = 1374.25
$ws.Range("L58").Value = 4381.6216
$ws.Range("M58").Value = -1171.25
$ws.Range("N58").Value = -4787.6216
# Row 62
$ws.Range("H62").Value = 3033.8235
$ws.Range("I62").Value = 2761.5
$ws.Range("J62").Value = 3422.8572
$ws.Range("K62").Value = 2761.5
$ws.Range("L62").Value = 3422.8572
$ws.Range("M62").Value = -2137.5
$ws.Range("N62").Value = -4670.8572
# Row 65
$ws.Range("H65").Value = 3033.8235
$ws.Range("I65").Value = 2761.5
$ws.Range("J65").Value = 3422.8572
$ws.Range("K65").Value = 13807.5
$ws.Range("L65").Value = 17114.286
$ws.Range("M65").Value = -10687.5
$ws.Range("N65").Value = -23354.286
# Row 99
$ws.Range("H99").Value = 188149
$ws.Range("I99").Value = 72168.28999999999
$ws.Range("J99").Value = 1000014
$ws.Range("K99").Value = 72168.28999999999
$ws.Range("L99").Value = 1000014
$ws.Range("M99").Value = -70670.28999999999
$ws.Range("N99").Value = -1003010
# Row 126
$ws.Range("H126").Value = 188149
$ws.Range("I126").Value = 72168.28999999999
$ws.Range("J126").Value = 1000014
$ws.Range("K126").Value = 216504.87
$ws.Range("L126").Value = 3000042
$ws.Range("M126").Value = -214034.87
$ws.Range("N126").Value = -3004982
# Row 134
$ws.Range("H134").Value = 1808.3103
$ws.Range("I134").Value = 1646.3617
$ws.Range("K134").Value = 4939.0851
$ws.Range("M134").Value = -2404.0851
# Row 136
$ws.Range("H136").Value = 3198.3936
$ws.Range("I136").Value = 1374.25
$ws.Range("J136").Value = 4381.6216
$ws.Range("K136").Value = 4122.75
$ws.Range("L136").Value = 13144.8648
$ws.Range("M136").Value = -1572.75
$ws.Range("N136").Value = -18244.8648

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 106
$ws.Range("H106").Value = 4078.5
$ws.Range("J106").Value = 4078.5
$ws.Range("L106").Value = 12235.5
$ws.Range("N106").Value = -14127.5
# Row 129
$ws.Range("H129").Value = 3596.9565
$ws.Range("J129").Value = 4383.3887
$ws.Range("L129").Value = 13150.1661
$ws.Range("N129").Value = -23150.1661
# Row 131
$ws.Range("H131").Value = 2101.6296
$ws.Range("I131").Value = 15112.5
$ws.Range("J131").Value = 1425.7402
$ws.Range("K131").Value = 45337.5
$ws.Range("L131").Value = 4277.2206
$ws.Range("M131").Value = -40297.5
$ws.Range("N131").Value = -14357.2206
# Row 134
$ws.Range("H134").Value = 16669456
$ws.Range("I134").Value = 26317300
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 78951900
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -78946830
$ws.Range("N134").Value = -25140

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
# Row 80
$ws.Range("H80").Value = 85909.25
$ws.Range("I80").Value = 2490.5
$ws.Range("J80").Value = 503003
$ws.Range("K80").Value = 2490.5
$ws.Range("L80").Value = 503003
$ws.Range("M80").Value = -1492.5
$ws.Range("N80").Value = -504999
# Row 83
$ws.Range("H83").Value = 85909.25
$ws.Range("I83").Value = 2490.5
$ws.Range("J83").Value = 503003
$ws.Range("K83").Value = 12452.5
$ws.Range("L83").Value = 2515015
$ws.Range("M83").Value = -7460.5
$ws.Range("N83").Value = -2524999
# Row 92
$ws.Range("H92").Value = 22000
$ws.Range("J92").Value = 22000
$ws.Range("L92").Value = 22000
$ws.Range("N92").Value = -25744
# Row 132
$ws.Range("H132").Value = 66669080
$ws.Range("I132").Value = 100001670
$ws.Range("J132").Value = 3903.6
$ws.Range("K132").Value = 300005010
$ws.Range("L132").Value = 11710.8
$ws.Range("M132").Value = -300002480
$ws.Range("N132").Value = -16770.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2445.2222
$ws.Range("I16").Value = 2567.8823
$ws.Range("J16").Value = 360
$ws.Range("K16").Value = 2567.8823
$ws.Range("L16").Value = 360
$ws.Range("M16").Value = -2397.8823
$ws.Range("N16").Value = -700

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4100.0713
$ws.Range("I62").Value = 4640.2
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 4640.2
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -4016.2
$ws.Range("N62").Value = -5048
# Row 65
$ws.Range("H65").Value = 4100.0713
$ws.Range("I65").Value = 4640.2
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 23201
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -20081
$ws.Range("N65").Value = -25240
# Row 136
$ws.Range("H136").Value = 4090008.2
$ws.Range("I136").Value = 10390
$ws.Range("J136").Value = 15875572
$ws.Range("K136").Value = 31170
$ws.Range("L136").Value = 47626716
$ws.Range("M136").Value = -28620
$ws.Range("N136").Value = -47631816
